$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Names" values for S3 and S4 from "Hoàng" to "Minh"
$ws.Range("S3").Value = "Minh"
$ws.Range("S4").Value = "Minh"

# Update the active selection to S9
$ws.Range("S9").Select()
